$d = $word.ActiveDocument

# The edit merges four "<id>...</id>" sequences (each previously split
# across three runs: the "<id>" open-tag run, a plain id-value run, and
# the "</id>" close-tag run) into a single run per occurrence. Using
# Find/Replace on the exact visible text collapses the found range into
# one run that inherits the formatting of the first character found
# (the "<id>" run's Courier-New / color 7f6000 formatting), which is
# exactly the formatting the merged run should end up with.

$targets = @(
    "<id>p077v_3</id>",
    "<id>p078r_1</id>",
    "<id>p078r_2</id>",
    "<id>p078r_3</id>"
)

foreach ($t in $targets) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $ok = $rng.Find.Execute($t, $false, $false, $false, $false, $false, $true, 1, $false, $t, 2)
    Write-Output "Replaced '$t': $ok"
}
